# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-12 (serial 45181) to 2023-09-13 (serial 45182).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
